$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.757.35'
$ws.Range("E2").Value = '  +0.08%  '

$ws.Range("D3").Value = '3.365.04'
$ws.Range("E3").Value = '  -0.47%  '

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.57'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +0.27%  '

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.52'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -1.96%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  -0.62%  '

$ws.Range("E9").Value = '  +2.95%  '

$ws.Range("E10").Value = '  -1.82%  '

$ws.Range("E11").Value = '  -4.01%  '

$ws.Range("D12").Value = '3.939.76'
$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("E13").Value = '  +0.45%  '

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.85'
$ws.Range("D14").Style = $origStyle

$ws.Range("D15").Value = '3.345.66'
$ws.Range("E15").Value = '  -0.91%  '

$ws.Range("D17").Value = '60.901.44'
$ws.Range("E17").Value = '  +0.15%  '

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.09'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  -1.87%  '

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.48'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  -3.38%  '

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.90'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  -0.47%  '

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '382.35'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -0.01%  '

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.22'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +2.20%  '

$ws.Range("E23").Value = '  -2.01%  '

$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("E25").Value = '  -5.70%  '

$ws.Range("E26").Value = '  +6.54%  '

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.12'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  -3.85%  '

$ws.Range("E29").Value = '  -1.43%  '

$ws.Range("E30").Value = '  -1.67%  '

$ws.Range("E31").Value = '  -0.03%  '

$ws.Range("E32").Value = '  -7.34%  '

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.93'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  -2.73%  '

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '167.15'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  +0.20%  '

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.80'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -1.71%  '

$ws.Range("E36").Value = '  -1.42%  '

$ws.Range("D37").Value = '3.399.97'
$ws.Range("E37").Value = '  -0.34%  '

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.43'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -3.16%  '

$ws.Range("E39").Value = '  -2.48%  '

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.39'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  -8.52%  '

$ws.Range("E41").Value = '  -1.20%  '

$ws.Range("E42").Value = '  -1.75%  '

$ws.Range("E43").Value = '  -2.73%  '

$ws.Range("E44").Value = '  -1.76%  '

$ws.Range("D45").Value = '2.439.57'
$ws.Range("E45").Value = '  -2.87%  '

$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("E47").Value = '  -2.87%  '

$ws.Range("E48").Value = '  -6.04%  '

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0257'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -4.87%  '

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.94'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  -6.12%  '

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.201'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  -2.81%  '
